$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = 44603
$ws.Range("I2").Value = 'Primera'
$ws.Range("J2").Value2 = 140
$ws.Range("K2").Value2 = 5500
$ws.Range("L2").Value2 = 6000
$ws.Range("M2").Value2 = 5750
$ws.Range("N2").Value = '$/caja 60 unidades'
$ws.Range("P2").Value2 = 96
$ws.Range("Q2").Value2 = 60

# Row 3
$ws.Range("D3").Value2 = 44785
$ws.Range("I3").Value = 'Primera'
$ws.Range("J3").Value2 = 130
$ws.Range("K3").Value2 = 7000
$ws.Range("L3").Value2 = 8000
$ws.Range("M3").Value2 = 7500
$ws.Range("N3").Value = '$/caja 60 unidades'
$ws.Range("P3").Value2 = 125
$ws.Range("Q3").Value2 = 60

# Row 4
$ws.Range("D4").Value2 = 45044
$ws.Range("I4").Value = 'Primera'
$ws.Range("J4").Value2 = 190
$ws.Range("K4").Value2 = 4000
$ws.Range("L4").Value2 = 5000
$ws.Range("M4").Value2 = 4526
$ws.Range("N4").Value = '$/caja 60 unidades'
$ws.Range("P4").Value2 = 75
$ws.Range("Q4").Value2 = 60

# Row 5
$ws.Range("D5").Value2 = 44494
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value2 = 120
$ws.Range("K5").Value2 = 5000
$ws.Range("L5").Value2 = 6000
$ws.Range("M5").Value2 = 5500
$ws.Range("N5").Value = '$/caja 60 unidades'
$ws.Range("P5").Value2 = 92
$ws.Range("Q5").Value2 = 60

# Row 6
$ws.Range("D6").Value2 = 44648
$ws.Range("I6").Value = 'Primera'
$ws.Range("J6").Value2 = 120
$ws.Range("K6").Value2 = 6500
$ws.Range("L6").Value2 = 7000
$ws.Range("M6").Value2 = 6750
$ws.Range("N6").Value = '$/caja 60 unidades'
$ws.Range("P6").Value2 = 112
$ws.Range("Q6").Value2 = 60

# Row 7
$ws.Range("D7").Value2 = 44589
$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value2 = 110
$ws.Range("K7").Value2 = 5000
$ws.Range("L7").Value2 = 6000
$ws.Range("M7").Value2 = 5500
$ws.Range("N7").Value = '$/caja 60 unidades'
$ws.Range("P7").Value2 = 92
$ws.Range("Q7").Value2 = 60

# Row 8
$ws.Range("D8").Value2 = 44421
$ws.Range("I8").Value = 'Primera'
$ws.Range("J8").Value2 = 100
$ws.Range("K8").Value2 = 8000
$ws.Range("L8").Value2 = 9000
$ws.Range("M8").Value2 = 8500
$ws.Range("N8").Value = '$/caja 60 unidades'
$ws.Range("P8").Value2 = 142
$ws.Range("Q8").Value2 = 60

# Row 9
$ws.Range("D9").Value2 = 44827
$ws.Range("I9").Value = 'Primera'
$ws.Range("J9").Value2 = 120
$ws.Range("K9").Value2 = 6000
$ws.Range("L9").Value2 = 7000
$ws.Range("M9").Value2 = 6500
$ws.Range("N9").Value = '$/caja 60 unidades'
$ws.Range("P9").Value2 = 108
$ws.Range("Q9").Value2 = 60

# Row 10
$ws.Range("D10").Value2 = 44740
$ws.Range("I10").Value = 'Primera'
$ws.Range("J10").Value2 = 120
$ws.Range("K10").Value2 = 6000
$ws.Range("L10").Value2 = 7000
$ws.Range("M10").Value2 = 6500
$ws.Range("N10").Value = '$/caja 60 unidades'
$ws.Range("P10").Value2 = 108
$ws.Range("Q10").Value2 = 60

# Row 11
$ws.Range("D11").Value2 = 44760
$ws.Range("I11").Value = 'Primera'
$ws.Range("J11").Value2 = 130
$ws.Range("K11").Value2 = 7000
$ws.Range("L11").Value2 = 7500
$ws.Range("M11").Value2 = 7250
$ws.Range("N11").Value = '$/caja 60 unidades'
$ws.Range("P11").Value2 = 121
$ws.Range("Q11").Value2 = 60

# Row 12
$ws.Range("D12").Value2 = 45079
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value2 = 130
$ws.Range("K12").Value2 = 4000
$ws.Range("L12").Value2 = 5000
$ws.Range("M12").Value2 = 4462
$ws.Range("N12").Value = '$/caja 60 unidades'
$ws.Range("P12").Value2 = 74
$ws.Range("Q12").Value2 = 60

# Row 13
$ws.Range("D13").Value2 = 44935
$ws.Range("I13").Value = 'Primera'
$ws.Range("J13").Value2 = 120
$ws.Range("K13").Value2 = 6000
$ws.Range("L13").Value2 = 7000
$ws.Range("M13").Value2 = 6500
$ws.Range("N13").Value = '$/caja 60 unidades'
$ws.Range("P13").Value2 = 108
$ws.Range("Q13").Value2 = 60

# Row 14
$ws.Range("D14").Value2 = 44362
$ws.Range("I14").Value = 'Primera'
$ws.Range("J14").Value2 = 120
$ws.Range("K14").Value2 = 8000
$ws.Range("L14").Value2 = 9000
$ws.Range("M14").Value2 = 8500
$ws.Range("N14").Value = '$/caja 60 unidades'
$ws.Range("P14").Value2 = 142
$ws.Range("Q14").Value2 = 60

# Row 15
$ws.Range("D15").Value2 = 44764
$ws.Range("I15").Value = 'Primera'
$ws.Range("J15").Value2 = 120
$ws.Range("K15").Value2 = 7000
$ws.Range("L15").Value2 = 8000
$ws.Range("M15").Value2 = 7500
$ws.Range("N15").Value = '$/caja 60 unidades'
$ws.Range("P15").Value2 = 125
$ws.Range("Q15").Value2 = 60

# Row 16
$ws.Range("D16").Value2 = 44963
$ws.Range("I16").Value = 'Primera'
$ws.Range("J16").Value2 = 130
$ws.Range("K16").Value2 = 4000
$ws.Range("L16").Value2 = 4500
$ws.Range("M16").Value2 = 4250
$ws.Range("N16").Value = '$/caja 60 unidades'
$ws.Range("P16").Value2 = 71
$ws.Range("Q16").Value2 = 60

# Row 17
$ws.Range("D17").Value2 = 44382
$ws.Range("I17").Value = 'Primera'
$ws.Range("J17").Value2 = 160
$ws.Range("K17").Value2 = 7000
$ws.Range("L17").Value2 = 8000
$ws.Range("M17").Value2 = 7438
$ws.Range("N17").Value = '$/caja 60 unidades'
$ws.Range("P17").Value2 = 124
$ws.Range("Q17").Value2 = 60

# Row 18
$ws.Range("D18").Value2 = 44242
$ws.Range("I18").Value = 'Primera'
$ws.Range("J18").Value2 = 160
$ws.Range("K18").Value2 = 5000
$ws.Range("L18").Value2 = 5500
$ws.Range("M18").Value2 = 5250
$ws.Range("N18").Value = '$/caja 60 unidades'
$ws.Range("P18").Value2 = 88
$ws.Range("Q18").Value2 = 60

# Row 19
$ws.Range("D19").Value2 = 44657
$ws.Range("I19").Value = 'Primera'
$ws.Range("J19").Value2 = 100
$ws.Range("K19").Value2 = 5000
$ws.Range("L19").Value2 = 5500
$ws.Range("M19").Value2 = 5250
$ws.Range("N19").Value = '$/caja 60 unidades'
$ws.Range("P19").Value2 = 88
$ws.Range("Q19").Value2 = 60

# Row 20
$ws.Range("D20").Value2 = 44281
$ws.Range("I20").Value = 'Primera'
$ws.Range("J20").Value2 = 120
$ws.Range("K20").Value2 = 5500
$ws.Range("L20").Value2 = 6000
$ws.Range("M20").Value2 = 5750
$ws.Range("N20").Value = '$/caja 60 unidades'
$ws.Range("P20").Value2 = 96
$ws.Range("Q20").Value2 = 60

# Row 21
$ws.Range("D21").Value2 = 44676
$ws.Range("I21").Value = 'Primera'
$ws.Range("J21").Value2 = 120
$ws.Range("K21").Value2 = 4000
$ws.Range("L21").Value2 = 4500
$ws.Range("M21").Value2 = 4250
$ws.Range("N21").Value = '$/caja 60 unidades'
$ws.Range("P21").Value2 = 71
$ws.Range("Q21").Value2 = 60

# Row 22
$ws.Range("D22").Value2 = 44669
$ws.Range("I22").Value = 'Primera'
$ws.Range("J22").Value2 = 130
$ws.Range("K22").Value2 = 4500
$ws.Range("L22").Value2 = 5000
$ws.Range("M22").Value2 = 4750
$ws.Range("N22").Value = '$/caja 60 unidades'
$ws.Range("P22").Value2 = 79
$ws.Range("Q22").Value2 = 60

# Row 23
$ws.Range("D23").Value2 = 44967
$ws.Range("I23").Value = 'Segunda'
$ws.Range("J23").Value2 = 50
$ws.Range("K23").Value2 = 4500
$ws.Range("L23").Value2 = 5000
$ws.Range("M23").Value2 = 4850
$ws.Range("N23").Value = '$/caja 90 unidades'
$ws.Range("P23").Value2 = 54
$ws.Range("Q23").Value2 = 90

# Row 24
$ws.Range("D24").Value2 = 44400
$ws.Range("I24").Value = 'Primera'
$ws.Range("J24").Value2 = 120
$ws.Range("K24").Value2 = 9000
$ws.Range("L24").Value2 = 10000
$ws.Range("M24").Value2 = 9500
$ws.Range("N24").Value = '$/caja 60 unidades'
$ws.Range("P24").Value2 = 158
$ws.Range("Q24").Value2 = 60

# Row 25
$ws.Range("D25").Value2 = 44627
$ws.Range("I25").Value = 'Primera'
$ws.Range("J25").Value2 = 120
$ws.Range("K25").Value2 = 4000
$ws.Range("L25").Value2 = 4500
$ws.Range("M25").Value2 = 4250
$ws.Range("N25").Value = '$/caja 60 unidades'
$ws.Range("P25").Value2 = 71
$ws.Range("Q25").Value2 = 60
